# regen sval data to filter save games
# Refresh the computed TB/d2S/K/IP/sum values (columns B,C,D,E,G) for rows 2-18.
# Column A (date) and column F (Win) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.1190320826869504, 0.306821227259698,  3.537761648806719,  0.4942365360607697, 4.457851494814137)
    3  = @(3.286832544864788,  1.655778082260271,  3.537761648806719,  0.4942365360607697, 8.974608811992548)
    4  = @(3.286832544864788,  1.655778082260271,  3.537761648806719,  0.4942365360607697, 8.974608811992548)
    5  = @(3.286832544864788,  1.655778082260271,  0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    6  = @(3.286832544864788,  1.655778082260271,  3.537761648806719,  0.4942365360607697, 8.974608811992548)
    7  = @(3.286832544864788,  1.655778082260271,  0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    8  = @(0.1190320826869504, 0.306821227259698,  0.1494219747398047, 0.4942365360607697, 1.069511820747223)
    9  = @(1.455362044514542,  1.655778082260271,  0.7527432677738641, 0.4942365360607697, 4.358119930609447)
    10 = @(3.286832544864788,  1.655778082260271,  3.537761648806719,  0.4942365360607697, 8.974608811992548)
    11 = @(3.286832544864788,  1.655778082260271,  3.537761648806719,  0.4942365360607697, 8.974608811992548)
    12 = @(1.455362044514542,  0.306821227259698,  22.3905356188092,   10.19245300693656,  34.34517189751999)
    13 = @(3.286832544864788,  1.655778082260271,  0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    14 = @(0.1190320826869504, 0.306821227259698,  3.537761648806719,  0.4942365360607697, 4.457851494814137)
    15 = @(0.04271373187048222,0.306821227259698,  0.7527432677738641, 0.4942365360607697, 1.596514762964814)
    16 = @(3.286832544864788,  1.655778082260271,  0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    17 = @(0.1190320826869504, 0.04071648406533734,0.1494219747398047, 0.4942365360607697, 0.8034070775528621)
    18 = @(0.6606524410359556, 1.655778082260271,  0.7527432677738641, 0.4942365360607697, 3.56341032713086)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Range("B$r").Value = $vals[0]
    $ws.Range("C$r").Value = $vals[1]
    $ws.Range("D$r").Value = $vals[2]
    $ws.Range("E$r").Value = $vals[3]
    $ws.Range("G$r").Value = $vals[4]
}
